# Updated cryptos list on Tue Feb 20 22:53:54 UTC 2024 with GitHub Actions
#
# Refreshes the Price (column D) and Volume(1h) (column E) figures for every
# coin row, and re-orders a handful of coins whose ranking changed (the
# Coin name + Link in column B/C move together with their row).
#
# NumberFormat is forced to Text ("@") before writing any Price value that
# would otherwise be auto-parsed by Excel as a number (losing formatting
# such as trailing zeros, e.g. "0.930" or "2.40").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "52.306.46"
$ws.Range("E2").Value = "  +0.62%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "2.997.84"
$ws.Range("E3").Value = "  +1.13%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.02%  "

# Row 5 - BNB
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "354.41"
$ws.Range("E5").Value = "  +0.23%  "

# Row 6 - Solana
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "108.58"
$ws.Range("E6").Value = "  -3.34%  "

# Row 7 - XRP
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.564"
$ws.Range("E7").Value = "  -0.31%  "

# Row 8 - USDC
$ws.Range("E8").Value = "  -0.05%  "

# Row 9 - Cardano
$ws.Range("E9").Value = "  -1.03%  "

# Row 10 - Avalanche
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "38.54"
$ws.Range("E10").Value = "  -3.02%  "

# Row 11 - TRON
$ws.Range("E11").Value = "  +2.42%  "

# Row 12 - Dogecoin
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0862"
$ws.Range("E12").Value = "  -4.11%  "

# Row 13 - Chainlink
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "19.36"
$ws.Range("E13").Value = "  -3.11%  "

# Row 14 - WrappedliquidstakedEther2.0
$ws.Range("D14").Value = "3.466.56"
$ws.Range("E14").Value = "  +0.85%  "

# Row 15 - Polkadot
$ws.Range("E15").Value = "  -3.07%  "

# Row 16 - WrappedEther
$ws.Range("D16").Value = "2.991.51"
$ws.Range("E16").Value = "  +0.82%  "

# Row 17 - Polygon
$ws.Range("E17").Value = "  +2.85%  "

# Row 18 - WrappedBTC
$ws.Range("D18").Value = "52.337.48"
$ws.Range("E18").Value = "  +0.45%  "

# Row 19 - ImmutableX
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.53"
$ws.Range("E19").Value = "  +6.03%  "

# Row 20 - Uniswap
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.57"
$ws.Range("E20").Value = "  -2.25%  "

# Row 21 - InternetComputer(DFINITY)
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.67"
$ws.Range("E21").Value = "  -6.10%  "

# Row 22 - ShibaInu
$ws.Range("D22").Value = "0.0₃0977"
$ws.Range("E22").Value = "  -1.56%  "

# Row 23 - Litecoin
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "69.63"
$ws.Range("E23").Value = "  -2.62%  "

# Row 24 - BitcoinCash
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "264.73"
$ws.Range("E24").Value = "  -2.39%  "

# Row 25 - PancakeSwap
$ws.Range("E25").Value = "  -2.44%  "

# Row 26 - Kaspa
$ws.Range("E26").Value = "  -0.79%  "

# Row 27 - EthereumClassic
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "26.98"
$ws.Range("E27").Value = "  -1.87%  "

# Row 28 - Filecoin
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.62"
$ws.Range("E28").Value = "  +0.33%  "

# Row 29 - Dai
$ws.Range("E29").Value = "  +0.08%  "

# Row 30 - Hedera
$ws.Range("E30").Value = "  -1.34%  "

# Row 31 - was RenderToken, now Cosmos (rows 31/32 swap rank)
$ws.Range("B31").Value = "Cosmos"
$ws.Range("C31").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "10.37"
$ws.Range("E31").Value = "  -3.63%  "

# Row 32 - was Cosmos, now RenderToken
$ws.Range("B32").Value = "RenderToken"
$ws.Range("C32").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.39"
$ws.Range("E32").Value = "  +2.83%  "

# Row 33 - InjectiveProtocol
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "36.67"
$ws.Range("E33").Value = "  -3.08%  "

# Row 34 - Toncoin
$ws.Range("E34").Value = "  +11.48%  "

# Row 35 - OKB
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "50.95"
$ws.Range("E35").Value = "  -4.28%  "

# Row 36 - VeChain
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0446"
$ws.Range("E36").Value = "  -1.34%  "

# Row 37 - FirstDigitalUSD
$ws.Range("E37").Value = "  -0.06%  "

# Row 38 - LidoDAOToken
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.21"
$ws.Range("E38").Value = "  -6.62%  "

# Row 39 - Celestia
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "17.98"
$ws.Range("E39").Value = "  -5.46%  "

# Row 40 - ARBITRUM
$ws.Range("E40").Value = "  -4.38%  "

# Row 41 - Stacks
$ws.Range("E41").Value = "  +0.36%  "

# Row 42 - Stellar
$ws.Range("E42").Value = "  -0.03%  "

# Row 43 - was EnergySwap, now Monero (rows 43/44 swap rank)
$ws.Range("B43").Value = "Monero"
$ws.Range("C43").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "124.65"
$ws.Range("E43").Value = "  +9.36%  "

# Row 44 - was Monero, now EnergySwap
$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "22.87"
$ws.Range("E44").Value = "  -4.61%  "

# Row 45 - WEMIXToken
$ws.Range("E45").Value = "  -1.42%  "

# Row 46 - Maker
$ws.Range("D46").Value = "2.126.40"
$ws.Range("E46").Value = "  -2.59%  "

# Row 47 - NEARProtocol
$ws.Range("E47").Value = "  -5.46%  "

# Row 48 - ApeXProtocol
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.40"
$ws.Range("E48").Value = "  -5.36%  "

# Row 49 - TheGraph
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.248"
$ws.Range("E49").Value = "  +1.95%  "

# Row 50 - was RocketPoolETH, now BEAM (RocketPoolETH drops out of the list)
$ws.Range("B50").Value = "BEAM"
$ws.Range("C50").Value = "https://coinranking.com/coin/cYYMfXF4u+beam-beam"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0333"
$ws.Range("E50").Value = "  -1.91%  "

# Row 51 - was BEAM, now SEI (new entry into the list)
$ws.Range("B51").Value = "SEI"
$ws.Range("C51").Value = "https://coinranking.com/coin/8nxCqs-uj+sei-sei"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.930"
$ws.Range("E51").Value = "  -1.04%  "
